{"js": "// Update the date line and every \"NN\u00f7N=\" division prompt in the table,\n// in document order. Each old value is unique in the document at the\n// moment it is searched (verified against the target diff), so a\n// straightforward ordered sequence of exact, case-sensitive searches\n// and in-place text replacements reproduces the authored edit while\n// preserving each run's existing formatting (font, size, etc.).\nconst replacements = [\n  [\"2025-09-10 Wednesday\", \"2025-09-11 Thursday\"],\n  [\"66\u00f75=\", \"35\u00f75=\"],\n  [\"62\u00f73=\", \"45\u00f74=\"],\n  [\"57\u00f73=\", \"81\u00f77=\"],\n  [\"87\u00f73=\", \"98\u00f79=\"],\n  [\"23\u00f76=\", \"72\u00f79=\"],\n  [\"66\u00f73=\", \"84\u00f77=\"],\n  [\"50\u00f74=\", \"59\u00f76=\"],\n  [\"43\u00f73=\", \"80\u00f78=\"],\n  [\"43\u00f79=\", \"89\u00f73=\"],\n  [\"33\u00f73=\", \"44\u00f73=\"],\n  [\"15\u00f75=\", \"46\u00f72=\"],\n  [\"99\u00f72=\", \"77\u00f75=\"],\n  [\"78\u00f72=\", \"21\u00f78=\"],\n  [\"66\u00f78=\", \"21\u00f72=\"],\n  [\"88\u00f74=\", \"78\u00f73=\"],\n  [\"96\u00f72=\", \"23\u00f72=\"],\n  [\"28\u00f78=\", \"94\u00f74=\"],\n  [\"73\u00f77=\", \"82\u00f74=\"],\n  [\"81\u00f73=\", \"91\u00f78=\"],\n  [\"75\u00f73=\", \"83\u00f73=\"],\n  [\"33\u00f77=\", \"29\u00f78=\"],\n  [\"54\u00f79=\", \"47\u00f72=\"],\n  [\"69\u00f76=\", \"12\u00f76=\"],\n  [\"47\u00f76=\", \"56\u00f79=\"],\n  [\"52\u00f76=\", \"54\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"NN\u00f7N=\" division prompt in the table,\n# in document order. Each old value is unique in the document at the\n# moment it is searched (verified against the target diff), so a plain\n# Find/Replace (wdReplaceAll, which here only ever touches the single\n# existing match) reproduces the authored edit while leaving each run's\n# existing formatting (font, size, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-10 Wednesday\", \"2025-09-11 Thursday\"),\n    @(\"66\u00f75=\", \"35\u00f75=\"),\n    @(\"62\u00f73=\", \"45\u00f74=\"),\n    @(\"57\u00f73=\", \"81\u00f77=\"),\n    @(\"87\u00f73=\", \"98\u00f79=\"),\n    @(\"23\u00f76=\", \"72\u00f79=\"),\n    @(\"66\u00f73=\", \"84\u00f77=\"),\n    @(\"50\u00f74=\", \"59\u00f76=\"),\n    @(\"43\u00f73=\", \"80\u00f78=\"),\n    @(\"43\u00f79=\", \"89\u00f73=\"),\n    @(\"33\u00f73=\", \"44\u00f73=\"),\n    @(\"15\u00f75=\", \"46\u00f72=\"),\n    @(\"99\u00f72=\", \"77\u00f75=\"),\n    @(\"78\u00f72=\", \"21\u00f78=\"),\n    @(\"66\u00f78=\", \"21\u00f72=\"),\n    @(\"88\u00f74=\", \"78\u00f73=\"),\n    @(\"96\u00f72=\", \"23\u00f72=\"),\n    @(\"28\u00f78=\", \"94\u00f74=\"),\n    @(\"73\u00f77=\", \"82\u00f74=\"),\n    @(\"81\u00f73=\", \"91\u00f78=\"),\n    @(\"75\u00f73=\", \"83\u00f73=\"),\n    @(\"33\u00f77=\", \"29\u00f78=\"),\n    @(\"54\u00f79=\", \"47\u00f72=\"),\n    @(\"69\u00f76=\", \"12\u00f76=\"),\n    @(\"47\u00f76=\", \"56\u00f79=\"),\n    @(\"52\u00f76=\", \"54\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
